$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update date in A1 (one month later: 45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update prices in column D for rows 31-34
$ws.Range("D31").Value = 6108
$ws.Range("D32").Value = 6845
$ws.Range("D33").Value = 7687
$ws.Range("D34").Value = 8424

# Refresh the merged-cell regions (unmerge/remerge) so the workbook's
# internal merge-cell bookkeeping matches the structure produced after
# the edits above.
$mergedRanges = @("A1:D1","B30:C30","A27:D27","A9:D9","B33:C33","B34:C34","B32:C32","B31:C31","A10:D10")
foreach ($r in $mergedRanges) {
    $ws.Range($r).UnMerge()
}
foreach ($r in $mergedRanges) {
    $ws.Range($r).Merge()
}
